# Deploying to gh-pages from @ NIH-NCPI/ncpi-fhir-ig-2@b701e861ff4aea87f49ab6a6b6da8d47ed8dfde7
# Update StructureDefinition-family-role workbook metadata + element text.

$wb = $excel.ActiveWorkbook

# --- Metadata sheet -------------------------------------------------------
$wsMeta = $wb.Worksheets.Item("Metadata")

# Date: regenerated IG build timestamp
$wsMeta.Range("B8").Value = "2025-06-13T15:45:04+00:00"

# FHIR Version: corrected from 4.3.0 (R4B) to 4.0.1 (R4)
$wsMeta.Range("B15").Value = "4.0.1"

# --- Elements sheet --------------------------------------------------------
$wsElem = $wb.Worksheets.Item("Elements")

# Extension row: ele-1 invariant text corrected (drop the Parameters-resource
# carve-out that doesn't apply to R4)
$wsElem.Range("AJ2").Value = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}`next-1:Must have either extensions or value[x], not both {extension.exists() != value.exists()}"

# Extension.id row: Type(s) corrected from "id" to "string"
$wsElem.Range("K3").Value = "string`n"

# Extension.value[x] row: Definition link corrected from R4B to R4
$wsElem.Range("M6").Value = "Value of extension - must be one of a constrained set of the data types (see [Extensibility](http://hl7.org/fhir/R4/extensibility.html) for a list)."
